$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (G=5503)
$ws.Range("H5").Value = 158
$ws.Range("I5").Value = 165.25
$ws.Range("K5").Value = 165.25
$ws.Range("M5").Value = -50.25

# Row 74 (G=5507)
$ws.Range("H74").Value = 18849.3
$ws.Range("I74").Value = 14748.833
$ws.Range("K74").Value = 14748.833
$ws.Range("M74").Value = -13812.833

# Row 77 (G=5507)
$ws.Range("H77").Value = 18849.3
$ws.Range("I77").Value = 14748.833
$ws.Range("K77").Value = 73744.16500000001
$ws.Range("M77").Value = -69064.16500000001

# Row 87 (G=10651)
$ws.Range("H87").Value = 89250
$ws.Range("J87").Value = 89250
$ws.Range("L87").Value = 89250
$ws.Range("N87").Value = -91746

# Row 90 (G=10651)
$ws.Range("H90").Value = 89250
$ws.Range("J90").Value = 89250
$ws.Range("L90").Value = 267750
$ws.Range("N90").Value = -280230

# Row 107 (G=27766)
$ws.Range("H107").Value = 223.82608
$ws.Range("I107").Value = 138.76471
$ws.Range("K107").Value = 138.76471
$ws.Range("M107").Value = 1781.23529

# Row 138 (G=44169)
$ws.Range("H138").Value = 3882.2
$ws.Range("I138").Value = 2994.5
$ws.Range("J138").Value = 4018.7693
$ws.Range("K138").Value = 8983.5
$ws.Range("L138").Value = 12056.3079
$ws.Range("M138").Value = -3843.5
$ws.Range("N138").Value = -22336.3079

$ws = $wb.Worksheets.Item("ARM")
# Row 16 (G=3775)
$ws.Range("H16").Value = 189.5
$ws.Range("J16").Value = 189.5
$ws.Range("L16").Value = 189.5
$ws.Range("N16").Value = -763.5

# Row 32 (G=44147)
$ws.Range("H32").Value = 4018809.2
$ws.Range("I32").Value = 11548.667
$ws.Range("K32").Value = 11548.667
$ws.Range("M32").Value = -11261.667

# Row 45 (G=27714)
$ws.Range("H45").Value = 6000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -6754

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (G=12526)
$ws.Range("H86").Value = 7388.6665
$ws.Range("J86").Value = 7999.7144
$ws.Range("L86").Value = 7999.7144
$ws.Range("N86").Value = -10245.7144

# Row 89 (G=12526)
$ws.Range("H89").Value = 7388.6665
$ws.Range("J89").Value = 7999.7144
$ws.Range("L89").Value = 39998.572
$ws.Range("N89").Value = -51230.572

$ws = $wb.Worksheets.Item("CRP")
# Row 17 (G=1823)
$ws.Range("H17").Value = 8332.666999999999
$ws.Range("I17").Value = 5000
$ws.Range("K17").Value = 5000
$ws.Range("M17").Value = -4826

# Row 25 (G=1895)
$ws.Range("H25").Value = 9999
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 31 (G=44023)
$ws.Range("H31").Value = 9760.652
$ws.Range("I31").Value = 9388.888999999999
$ws.Range("K31").Value = 9388.888999999999
$ws.Range("M31").Value = -9093.888999999999

# Row 34 (G=44023)
$ws.Range("H34").Value = 9760.652
$ws.Range("I34").Value = 9388.888999999999
$ws.Range("K34").Value = 9388.888999999999
$ws.Range("M34").Value = -9186.888999999999

# Row 59 (G=1942)
$ws.Range("H59").Value = 59791.9
$ws.Range("J59").Value = 70001.75
$ws.Range("L59").Value = 70001.75
$ws.Range("N59").Value = -72291.75

# Row 74 (G=10636)
$ws.Range("H74").Value = 80000
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81748

# Row 77 (G=10636)
$ws.Range("H77").Value = 80000
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -248736

# Row 134 (G=44020)
$ws.Range("H134").Value = 3868.5
$ws.Range("I134").Value = 2880.1428
$ws.Range("K134").Value = 8640.428400000001
$ws.Range("M134").Value = -6105.428400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 104 (G=19807)
$ws.Range("H104").Value = 12737.5
$ws.Range("I104").Value = 7300
$ws.Range("K104").Value = 21900
$ws.Range("M104").Value = -19279

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (G=1681)
$ws.Range("H5").Value = 10003333
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 48 (G=4337)
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970

# Row 80 (G=12521)
$ws.Range("H80").Value = 3958.3333
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996

# Row 83 (G=12521)
$ws.Range("H83").Value = 3958.3333
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984

# Row 126 (G=36184)
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9530
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G=5277)
$ws.Range("H22").Value = 738.6667
$ws.Range("I22").Value = 683.5
$ws.Range("K22").Value = 683.5
$ws.Range("M22").Value = -388.5

# Row 27 (G=5277)
$ws.Range("H27").Value = 738.6667
$ws.Range("I27").Value = 683.5
$ws.Range("K27").Value = 683.5
$ws.Range("M27").Value = -576.5

# Row 39 (G=1708)
$ws.Range("H39").Value = 8000
$ws.Range("I39").Value = 8000
$ws.Range("K39").Value = 8000
$ws.Range("M39").Value = -7540

# Row 40 (G=36248)
$ws.Range("H40").Value = 5510.625
$ws.Range("I40").Value = 4417.727
$ws.Range("K40").Value = 4417.727
$ws.Range("M40").Value = -4281.727

# Row 46 (G=5282)
$ws.Range("H46").Value = 7081.25
$ws.Range("I46").Value = 1650
$ws.Range("K46").Value = 1650
$ws.Range("M46").Value = -1462

# Row 61 (G=27740)
$ws.Range("H61").Value = 5349.6113
$ws.Range("I61").Value = 4482.75
$ws.Range("K61").Value = 4482.75
$ws.Range("M61").Value = -4280.75

# Row 93 (G=19993)
$ws.Range("H93").Value = 1491.3182
$ws.Range("I93").Value = 1599.2
$ws.Range("J93").Value = 1459.5883
$ws.Range("K93").Value = 1599.2
$ws.Range("L93").Value = 1459.5883
$ws.Range("M93").Value = -351.2
$ws.Range("N93").Value = -3955.5883

# Row 113 (G=27740)
$ws.Range("H113").Value = 5349.6113
$ws.Range("I113").Value = 4482.75
$ws.Range("K113").Value = 4482.75
$ws.Range("M113").Value = -2312.75

# Row 136 (G=44060)
$ws.Range("H136").Value = 4620.1
$ws.Range("I136").Value = 3837.75
$ws.Range("J136").Value = 7749.5
$ws.Range("K136").Value = 11513.25
$ws.Range("L136").Value = 23248.5
$ws.Range("M136").Value = -8963.25
$ws.Range("N136").Value = -28348.5

$ws = $wb.Worksheets.Item("WVR")
# Row 2 (G=3307)
$ws.Range("H2").Value = 113611.11
$ws.Range("I2").Value = 144642.86
$ws.Range("K2").Value = 144642.86
$ws.Range("M2").Value = -144530.86

# Row 15 (G=2670)
$ws.Range("H15").Value = 69999
$ws.Range("J15").Value = 69999
$ws.Range("L15").Value = 69999
$ws.Range("N15").Value = -70575

# Row 23 (G=3325)
$ws.Range("H23").Value = 2497.5
$ws.Range("I23").Value = 2495
$ws.Range("K23").Value = 2495
$ws.Range("M23").Value = -2266

# Row 54 (G=3413)
$ws.Range("H54").Value = 62706.855
$ws.Range("J54").Value = 62706.855
$ws.Range("L54").Value = 62706.855
$ws.Range("N54").Value = -63746.855

# Row 127 (G=35414)
$ws.Range("H127").Value = 69998
$ws.Range("I127").Value = 69998
$ws.Range("K127").Value = 69998
$ws.Range("M127").Value = -65038

# Row 132 (G=44029)
$ws.Range("H132").Value = 2049.2856
$ws.Range("I132").Value = 2049.2856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6147.8568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3617.8568
$ws.Range("N132").ClearContents()

# Row 136 (G=44031)
$ws.Range("H136").Value = 4504.8887
$ws.Range("I136").Value = 2936.2856
$ws.Range("J136").Value = 9995
$ws.Range("K136").Value = 8808.856800000001
$ws.Range("L136").Value = 29985
$ws.Range("M136").Value = -6258.856800000001
$ws.Range("N136").Value = -35085
